$wb = $excel.ActiveWorkbook

# --- Sheet: Risk Aware Portfolio ---
$ws1 = $wb.Worksheets.Item("Risk Aware Portfolio")
$ws1.Range("B2").Value = 0.2061
$ws1.Range("C2").Value = 103030
$ws1.Range("D2").Value = 2870

$ws1.Range("B3").Value = 0.1843
$ws1.Range("C3").Value = 92160
$ws1.Range("D3").Value = 8014

$ws1.Range("B4").Value = 0.2024
$ws1.Range("C4").Value = 101210
$ws1.Range("D4").Value = 1260

$ws1.Range("B5").Value = 0.1945
$ws1.Range("C5").Value = 97255
$ws1.Range("D5").Value = 2955

$ws1.Range("B6").Value = 0.2127
$ws1.Range("C6").Value = 106350
$ws1.Range("D6").Value = 106350

# --- Sheet: Risk Aware Performance ---
$ws2 = $wb.Worksheets.Item("Risk Aware Performance")
$ws2.Range("B2").Value = 0.06765256308685724
$ws2.Range("B3").Value = 0.3542435047398688
$ws2.Range("B4").Value = 0.1345192288616552

# --- Sheet: Max Sharpe Portfolio ---
$ws3 = $wb.Worksheets.Item("Max Sharpe Portfolio")
$ws3.Range("B2").Value = 0.3432
$ws3.Range("C2").Value = 171585
$ws3.Range("D2").Value = 4780

$ws3.Range("B3").Value = 0.0313
$ws3.Range("C3").Value = 15650
$ws3.Range("D3").Value = 1361

$ws3.Range("B4").Value = 0.2641
$ws3.Range("C4").Value = 132050
$ws3.Range("D4").Value = 1644

$ws3.Range("B5").Value = 0.204
$ws3.Range("C5").Value = 101980
$ws3.Range("D5").Value = 3099

$ws3.Range("B6").Value = 0.1575
$ws3.Range("C6").Value = 78735
$ws3.Range("D6").Value = 78735

# --- Sheet: Max Sharpe Performance ---
$ws4 = $wb.Worksheets.Item("Max Sharpe Performance")
$ws4.Range("B2").Value = 0.07857468323672587
$ws4.Range("B3").Value = 0.3949529727706366
$ws4.Range("B4").Value = 0.1483079943057987
